$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G: every row 2-51 changes from "4" to "5" (Text values, not numbers).
# A leading apostrophe forces Excel to treat the numeric-looking string as text,
# matching the original inline-string "Text" cell type; resetting the style back
# to "Normal" afterwards avoids leaving a stray quote-prefix number format behind.
$gRange = $ws.Range("G2:G51")
$gRange.Value = "'5"
$gRange.Style = "Normal"

# Column D: updated price readings for the rows whose price changed in this refresh.
$ws.Range("D2").Value = "'248.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("D4").Value = "'5.347"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05691"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.400"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.314"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8131"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.9114"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1403"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07400"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.03110"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03023"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09363"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.748"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001572"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04771"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Value = "'0.0005796"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.006469"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.004994"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.001026"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'0.0001501"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'3.698"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'2.147"
$ws.Range("D25").Style = "Normal"
$ws.Range("D40").Value = "'0.03965"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006854"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Value = "'0.002712"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.007456"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005897"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.5006"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").Value = "'0.01011"
$ws.Range("D50").Style = "Normal"
